# "Generate Report for Handback"
# The b5fa3b00-7fc8-45f1-91f6-e18accf757cd.md handback file has caught up to
# the latest handoff, so its status flips from "Ready for handoff" to
# "Handed back: in sync with en-US" everywhere it is reported, the
# Latest Handback DateTime is refreshed, and the stale "version mismatch"
# error is cleared now that the file is in sync.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for b5fa3b00-...md (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: row for b5fa3b00-...md (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-10-21 00:34:35"
$wsZhCn.Range("P3").Value = ""
$wsZhCn.Columns("P:P").AutoFit()

# --- de-de sheet: row for b5fa3b00-...md (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-10-21 00:34:53"
$wsDeDe.Range("P3").Value = ""
$wsDeDe.Columns("P:P").AutoFit()

Write-Output "Report generated for handback of b5fa3b00-7fc8-45f1-91f6-e18accf757cd.md"
